$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 221

# Add the new "Rank" header in column I, reusing the same formatting
# (bold, bordered, centered) as the other header cells by copying H1's format.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "Rank"

# Compute the rank of each student based on their FinalGPA (column G),
# highest GPA first (standard competition ranking, ties share a rank).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Formula = "=RANK(G$r,`$G`$2:`$G`$$lastRow,0)"
}
